$wb = $excel.ActiveWorkbook

# Column F ("想去人数") updates that apply identically to both the
# "展览" sheet and the "全部类型" sheet (they mirror the same data).
$updates = @{
    5  = 2723
    6  = 188
    7  = 142
    9  = 1442
    13 = 1220
    22 = 2653
    24 = 308
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
